# Fill in student IDs and names, and mark progress on several
# basic/extra feature rows (implemented recursive ray reflections, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# --- Team member identification (was "<student id>" / "<student name>" placeholders) ---
$ws.Range("D5").Value = 5748542
$ws.Range("E5").Value = 5702364
$ws.Range("F5").Value = 5482526

$ws.Range("D6").Value = "Vlad-Stefan Graure"
$ws.Range("E6").Value = "Rares Burghelea"
$ws.Range("F6").Value = "Ariel Potolski Eilat"

# --- Basic features progress (rows 8-15) ---
# row 8  generation and traversal of acceleration data-structure -> member2 100%
$ws.Range("E8").Value = 100
# row 9  implementation of shading models -> member3 100%
$ws.Range("F9").Value = 100
# row 10 implementation of recursive ray reflections -> member1 100%
$ws.Range("D10").Value = 100
# row 11 implementation of recursive ray transparency -> member3 100%
$ws.Range("F11").Value = 100
# row 12 normal interpolation with barycentric coordinates -> member3 100%
$ws.Range("F12").Value = 100
# row 13 implementation of texture mapping -> member3 100%
$ws.Range("F13").Value = 100
# row 14 implementation of lights and shadows -> member1 100%
$ws.Range("D14").Value = 100
# row 15 implementation of multisampling -> member3 100%
$ws.Range("F15").Value = 100

# --- Extra features progress (rows 19-24) ---
# row 19 Environment maps -> member3 100%
$ws.Range("F19").Value = 100
# row 20 SAH+binning as splitting criterion for BVH -> member3 100%
$ws.Range("F20").Value = 100
# row 21 Motion blur -> member2 100%
$ws.Range("E21").Value = 100
# row 22 Bloom filter -> member3 100%
$ws.Range("F22").Value = 100
# row 23 Glossy reflections -> member1 100%
$ws.Range("D23").Value = 100
# row 24 Depth of field -> member1 100%
$ws.Range("D24").Value = 100
